# Applies the commit "Fruta / hortaliza, semanal" update to the Tuna sheet.
# The underlying change inserts two new weekly observation rows (new rows 45
# and 46) into the "Femacal de La Calera - Tuna" data table, pushing the
# previously-existing rows 45-92 down to rows 47-94.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 45. This shifts old rows
# 45..92 down to become rows 47..94 (matching the target layout exactly).
$ws.Rows.Item(45).EntireRow.Insert()
$ws.Rows.Item(45).EntireRow.Insert()

# --- Fill the two newly inserted rows with the new weekly data ---
# These columns are constant for every data row in this sheet.
$marketId = 3
$market   = "Femacal de La Calera"
$region   = "Coquimbo"
$codreg   = 5
$tipo     = "Fruta"
$prodId   = 100107
$prod     = "Otros"
$catId    = 100107011
$cat      = "Tuna"
$variedad = "Sin especificar"

# New row 45
$ws.Cells.Item(45, 1).Value  = $marketId
$ws.Cells.Item(45, 2).Value  = $market
$ws.Cells.Item(45, 3).Value  = $region
$ws.Cells.Item(45, 4).Value  = 44589
$ws.Cells.Item(45, 5).Value  = $codreg
$ws.Cells.Item(45, 6).Value  = $tipo
$ws.Cells.Item(45, 7).Value  = $prodId
$ws.Cells.Item(45, 8).Value  = $prod
$ws.Cells.Item(45, 9).Value  = $catId
$ws.Cells.Item(45, 10).Value = $cat
$ws.Cells.Item(45, 11).Value = $variedad
$ws.Cells.Item(45, 12).Value = "Primera"
$ws.Cells.Item(45, 13).Value = 55
$ws.Cells.Item(45, 14).Value = 19000
$ws.Cells.Item(45, 15).Value = 19000
$ws.Cells.Item(45, 16).Value = 19000
$ws.Cells.Item(45, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(45, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(45, 19).Value = 950
$ws.Cells.Item(45, 20).Value = 20

# New row 46
$ws.Cells.Item(46, 1).Value  = $marketId
$ws.Cells.Item(46, 2).Value  = $market
$ws.Cells.Item(46, 3).Value  = $region
$ws.Cells.Item(46, 4).Value  = 44589
$ws.Cells.Item(46, 5).Value  = $codreg
$ws.Cells.Item(46, 6).Value  = $tipo
$ws.Cells.Item(46, 7).Value  = $prodId
$ws.Cells.Item(46, 8).Value  = $prod
$ws.Cells.Item(46, 9).Value  = $catId
$ws.Cells.Item(46, 10).Value = $cat
$ws.Cells.Item(46, 11).Value = $variedad
$ws.Cells.Item(46, 12).Value = "Segunda"
$ws.Cells.Item(46, 13).Value = 50
$ws.Cells.Item(46, 14).Value = 16000
$ws.Cells.Item(46, 15).Value = 16000
$ws.Cells.Item(46, 16).Value = 16000
$ws.Cells.Item(46, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(46, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(46, 19).Value = 800
$ws.Cells.Item(46, 20).Value = 20

Write-Output "done"
